$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Kris Demor"
$ws.Range("B3").Value = 40
